# Developed test methods for Filters
#
# This script applies the following changes to TestScript_3_0_customer.xlsx:
#  1. Row 15, column B ("Execute") value changes from "No" to "Yes".
#  2. A brand-new row 60 is appended describing the "Filters" test case,
#     mirroring the layout/format of the existing "Scan Code" style rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Flip the Execute flag for row 15 (SignUp with Add Address) to Yes
# ---------------------------------------------------------------------
$ws.Cells.Item(15, 2).Value = "Yes"

# ---------------------------------------------------------------------
# 2) Build new row 60 - "Verify Filters with Transactions Types,
#    Sub Types and Status"
# ---------------------------------------------------------------------

# Re-use the formatting of an existing, similarly laid out row (57) for
# columns A:I so the new row picks up the same (already existing) style
# indexes instead of Excel fabricating brand new ones.
$ws.Range("A57:I57").Copy()
$ws.Range("A60:I60").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# Column K / J formatting comes from other rows that already carry the
# blank, wrap-text style used at the right-hand edge of the table.
$ws.Range("K59").Copy()
$ws.Range("K60").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("J15").Copy()
$ws.Range("J60").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# Populate the new/unique strings first in the same order they show up
# in the regenerated shared-strings table (C, G, I, A), then backfill
# the cells that reuse already-existing shared strings.
$ws.Cells.Item(60, 3).Value = "testdata_3_0_customer.xls,filters"
$ws.Cells.Item(60, 7).Value = "Filters"
$ws.Cells.Item(60, 9).Value = "coyni_mobile.tests.DashBoardTest,`ntestFilters,`n-ptransactionHeading,`n-ptransDtlsHeading,`n-ptransactionType,`n-pfromAmount,`n-ptoAmount,`n-ptransactionType"
$ws.Cells.Item(60, 1).Value = "Verify Filters with Transactions Types,Sub Types and Status"

$ws.Cells.Item(60, 2).Value = "No"
$ws.Cells.Item(60, 4).Value = "RunRangeOfIterations"
$ws.Range("E60").Value = "'3"
$ws.Range("F60").Value = "'3"
$ws.Cells.Item(60, 8).Value = "coyni_mobile.tests.LoginTest,`ntestLogin,`n-pemail,`n-ppassword,`n-ppin"

$ws.Rows.Item(60).RowHeight = 120

# ---------------------------------------------------------------------
# 3) Restore the view/selection state as closely as possible
# ---------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 14
$ws.Range("B15").Select()
